$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns E:L one column to the right (F:M) to make room for the new
# "Empty Column 1" column at E. Work right-to-left so the single-column
# moves never overwrite data we still need to move.
$ws.Range("L1:L3").Cut($ws.Range("M1:M3"))
$ws.Range("K1:K3").Cut($ws.Range("L1:L3"))
$ws.Range("J1:J3").Cut($ws.Range("K1:K3"))
$ws.Range("I1:I3").Cut($ws.Range("J1:J3"))
$ws.Range("H1:H3").Cut($ws.Range("I1:I3"))
$ws.Range("G1:G3").Cut($ws.Range("H1:H3"))
$ws.Range("F1:F3").Cut($ws.Range("G1:G3"))
$ws.Range("E1:E3").Cut($ws.Range("F1:F3"))

# New, empty columns: E ("Empty Column 1") and N ("Empty Column 2") -
# only the header row carries a value, the data rows stay blank.
$ws.Range("E1").Value = "Empty Column 1"
$ws.Range("N1").Value = "Empty Column 2"

# The hyperlink that used to live on I2 is now on J2 (it rode along with
# the column shift above, but the hyperlink definition itself still
# points at the old address) - recreate it in the right place.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:nicht@meine.mail")
$ws.Range("J2").Style = "Link"

# Selection moved to E1 in the saved file.
$ws.Range("E1").Select()
